$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 82; this shifts the existing rows 82-189 down
# to 83-190 (and pulls formatting, e.g. the date style on column D, down
# with them automatically).
$ws.Rows("82:82").Insert()

# Populate the newly inserted row 82 with the new record's data.
$ws.Range("A82").Value = 10
$ws.Range("B82").Value = "Vega Modelo de Temuco"
$ws.Range("C82").Value = "La Araucanía"
$ws.Range("D82").Value = (Get-Date -Year 2022 -Month 3 -Day 18 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E82").Value = 9
$ws.Range("F82").Value = 100112005
$ws.Range("G82").Value = "Puerro"
$ws.Range("H82").Value = "Azul de Maquehue"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 30
$ws.Range("K82").Value = 10000
$ws.Range("L82").Value = 10000
$ws.Range("M82").Value = 10000
$ws.Range("N82").Value = "$/docena de paquetes"
$ws.Range("O82").Value = "Provincia de Cautín"
$ws.Range("P82").Value = 833
$ws.Range("Q82").Value = 12
$ws.Range("R82").Value = "Hortaliza"
